# Add a new entry (row 9) to the LeetCode tracking sheet for the
# "Regular Expression Matching" problem - naive implementation WIP.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the new row with its values -----------------------------
$ws.Range("A9").Value = "Hard"
$ws.Range("B9").Value = "regular expression match"
$ws.Range("C9").Value = "my_re_match.py"
$ws.Range("D9").Value = "N"
$ws.Range("E9").Value = "link"
$ws.Range("F9").Value = "string parsing"
$ws.Range("G9").Value = "TBD! (for now: try to find different cases / scenarios)"

# --- Hyperlink for the "link" cell, same pattern as the other rows ---
$ws.Hyperlinks.Add($ws.Range("E9"), "https://leetcode.com/problems/regular-expression-matching/", "", "", "link")

# --- Match the look & feel (centered / wrapped) of the other rows ----
$ws.Range("A9:G9").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A9:G9").VerticalAlignment = -4108     # xlCenter
$ws.Range("A9:G9").WrapText = $true

# --- Highlight the "Hard" level cell (red text on yellow fill) -------
$ws.Range("A9").Font.Color = 255       # red   (FFFF0000)
$ws.Range("A9").Interior.Color = 65535 # yellow (FFFFFF00)

# --- Row grew to two lines tall ---------------------------------------
$ws.Rows.Item(9).RowHeight = 28.8

# --- Leave the selection where the author ended up --------------------
$ws.Range("B10").Select() | Out-Null
